$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1. Copy the formatting (bold, centered,
# bordered header style) from the neighboring "sum" header (G1) so the new
# column matches the existing header row, then overwrite the copied text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for the new column in row 2 (plain,
# unstyled numeric cell like the other data cells).
$ws.Range("H2").Value = 0
